$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy the style (xlPasteFormats) of a source cell onto a destination cell,
# reusing existing style indices instead of synthesizing new ones.
function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Helper: copy the value (xlPasteValues) of a source cell onto a destination cell,
# so that string cells reuse the existing sharedStrings entry instead of duplicating it.
function Copy-CellValue($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4163) | Out-Null
}

# --- Row 19 ---------------------------------------------------------------
Copy-CellFormat "A18" "A19"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 18
Copy-CellValue "C15" "C19"
$ws.Range("D19").Value = 1
Copy-CellValue "E15" "E19"
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("J19").Value = 3
$ws.Range("K19").Value = 6
$ws.Range("L19").Value = 3
Copy-CellValue "M18" "M19"
Copy-CellValue "N18" "N19"
$ws.Range("O19").Value = "PmReviewed"
Copy-CellValue "P18" "P19"
Copy-CellValue "Q18" "Q19"
$ws.Range("S19").Value = 2000000
$ws.Range("U19").Value = 4

# --- Row 20 -----------------------------------------------------------------
Copy-CellFormat "A18" "A20"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 19
Copy-CellValue "C15" "C20"
$ws.Range("D20").Value = 1
Copy-CellValue "E15" "E20"
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 7
$ws.Range("L20").Value = 3
Copy-CellValue "M2" "M20"
Copy-CellValue "N3" "N20"
$ws.Range("O20").Value = "HrApproved"
Copy-CellValue "P18" "P20"
Copy-CellValue "Q18" "Q20"
$ws.Range("S20").Value = 2000000
$ws.Range("U20").Value = 4

# --- Row 21 -----------------------------------------------------------------
Copy-CellFormat "A18" "A21"
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 20
Copy-CellValue "C15" "C21"
$ws.Range("D21").Value = 1
Copy-CellValue "E15" "E21"
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("J21").Value = 3
$ws.Range("K21").Value = 12
$ws.Range("L21").Value = 3
Copy-CellValue "M18" "M21"
Copy-CellValue "N18" "N21"
$ws.Range("O21").Value = "ReOpen"
Copy-CellValue "P18" "P21"
Copy-CellValue "Q18" "Q21"
$ws.Range("S21").Value = 2000000
$ws.Range("U21").Value = 4

# --- Row 22 -----------------------------------------------------------------
Copy-CellFormat "A18" "A22"
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 21
Copy-CellValue "C15" "C22"
$ws.Range("D22").Value = 1
Copy-CellValue "E15" "E22"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 12
$ws.Range("L22").Value = 3
Copy-CellValue "M18" "M22"
Copy-CellValue "N18" "N22"
$ws.Range("O22").Value = "Rejected"
Copy-CellValue "P18" "P22"
Copy-CellValue "Q18" "Q22"
$ws.Range("S22").Value = 2000000
$ws.Range("U22").Value = 4

$ws.Application.CutCopyMode = $false

# Update the view's active selection to match the post-edit state.
$ws.Range("I24").Select() | Out-Null

Write-Host "edit complete"
